$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1952054794520548
$ws.Range("C2").Value = 0.541095890410959
$ws.Range("J2").Value = 0.01027397260273973
$ws.Range("P2").Value = 0.1438356164383562
$ws.Range("S2").Value = 0.1095890410958904
$ws.Range("B3").Value = 0.006097560975609756
$ws.Range("C3").Value = 0.03048780487804878
$ws.Range("J3").Value = 0.0426829268292683
$ws.Range("P3").Value = 0.7378048780487805
$ws.Range("S3").Value = 0.1829268292682927
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.5882352941176471
$ws.Range("S4").Value = 0.3529411764705883
$ws.Range("B6").Value = 0.04744525547445255
$ws.Range("D6").Value = 0.0218978102189781
$ws.Range("F6").Value = 0.0948905109489051
$ws.Range("J6").Value = 0.2664233576642336
$ws.Range("O6").Value = 0.01094890510948905
$ws.Range("Q6").Value = 0.0948905109489051
$ws.Range("R6").Value = 0.1094890510948905
$ws.Range("S6").Value = 0.354014598540146
$ws.Range("B7").Value = 0.08387096774193549
$ws.Range("D7").Value = 0.05161290322580645
$ws.Range("F7").Value = 0.05806451612903226
$ws.Range("J7").Value = 0.1354838709677419
$ws.Range("O7").Value = 0.03225806451612903
$ws.Range("Q7").Value = 0.1741935483870968
$ws.Range("R7").Value = 0.08387096774193549
$ws.Range("S7").Value = 0.3806451612903226
$ws.Range("B8").Value = 0.1094147582697201
$ws.Range("D8").Value = 0.02290076335877863
$ws.Range("F8").Value = 0.06361323155216285
$ws.Range("J8").Value = 0.1272264631043257
$ws.Range("O8").Value = 0.02544529262086514
$ws.Range("Q8").Value = 0.1653944020356234
$ws.Range("R8").Value = 0.089058524173028
$ws.Range("S8").Value = 0.3969465648854962
$ws.Range("B9").Value = 0.09166666666666666
$ws.Range("D9").Value = 0.01666666666666667
$ws.Range("E9").Value = 0.004166666666666667
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1291666666666667
$ws.Range("O9").Value = 0.01666666666666667
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.1041666666666667
$ws.Range("S9").Value = 0.4041666666666667
$ws.Range("B10").Value = 0.1044885945548197
$ws.Range("D10").Value = 0.01766004415011038
$ws.Range("E10").Value = 0.001471670345842531
$ws.Range("F10").Value = 0.08535688005886681
$ws.Range("J10").Value = 0.1066961000735835
$ws.Range("O10").Value = 0.01324503311258278
$ws.Range("Q10").Value = 0.2001471670345843
$ws.Range("R10").Value = 0.09050772626931568
$ws.Range("S10").Value = 0.3804267844002943
$ws.Range("G11").Value = 0.1149425287356322
$ws.Range("J11").Value = 0.09961685823754789
$ws.Range("K11").Value = 0.1647509578544061
$ws.Range("L11").Value = 0.6091954022988506
$ws.Range("S11").Value = 0.01149425287356322
$ws.Range("G12").Value = 0.70625
$ws.Range("J12").Value = 0.24375
$ws.Range("K12").Value = 0.00625
$ws.Range("L12").Value = 0.0125
$ws.Range("S12").Value = 0.03125
$ws.Range("F15").Value = 0.0135593220338983
$ws.Range("H15").Value = 0.176271186440678
$ws.Range("I15").Value = 0.09491525423728814
$ws.Range("J15").Value = 0.3864406779661017
$ws.Range("K15").Value = 0.06440677966101695
$ws.Range("M15").Value = 0.0135593220338983
$ws.Range("O15").Value = 0.1016949152542373
$ws.Range("S15").Value = 0.1491525423728814
$ws.Range("F16").Value = 0.01612903225806452
$ws.Range("H16").Value = 0.1989247311827957
$ws.Range("I16").Value = 0.08064516129032258
$ws.Range("J16").Value = 0.3279569892473118
$ws.Range("K16").Value = 0.1236559139784946
$ws.Range("M16").Value = 0.01612903225806452
$ws.Range("O16").Value = 0.1021505376344086
$ws.Range("S16").Value = 0.1344086021505376
$ws.Range("F17").Value = 0.02320185614849188
$ws.Range("H17").Value = 0.1392111368909513
$ws.Range("I17").Value = 0.1160092807424594
$ws.Range("J17").Value = 0.431554524361949
$ws.Range("K17").Value = 0.09280742459396751
$ws.Range("M17").Value = 0.006960556844547564
$ws.Range("N17").Value = 0.002320185614849188
$ws.Range("O17").Value = 0.0765661252900232
$ws.Range("S17").Value = 0.111368909512761
$ws.Range("F18").Value = 0.04888888888888889
$ws.Range("H18").Value = 0.1511111111111111
$ws.Range("I18").Value = 0.09777777777777778
$ws.Range("J18").Value = 0.4311111111111111
$ws.Range("K18").Value = 0.04888888888888889
$ws.Range("M18").Value = 0.01333333333333333
$ws.Range("O18").Value = 0.08
$ws.Range("S18").Value = 0.1288888888888889
$ws.Range("F19").Value = 0.0231814548361311
$ws.Range("H19").Value = 0.1710631494804157
$ws.Range("I19").Value = 0.1015187849720224
$ws.Range("J19").Value = 0.4044764188649081
$ws.Range("K19").Value = 0.09512390087929656
$ws.Range("M19").Value = 0.01678657074340528
$ws.Range("O19").Value = 0.09832134292565947
